$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 46 data (ajout ssr, had et psy)
$ws.Range("B46").Value = "c"
$ws.Range("D46").Value = 182
$ws.Range("F46").Value = "ZAD"

# Reflect the scrolled view / new selection on the sheet (topLeftCell A12, selection D47)
$ws.Range("D47").Select() | Out-Null
